# #5: insurance, claim, debt, investment done
# Rework the "保險" (insurance) sheet so that:
#  - row 1 becomes a proper header row (company, name, owner,
#    property_category, category, date, legislator_name, legislator_id,
#    source_file, index) instead of duplicating row 2's sample data,
#  - rows 2 & 3 gain the same trailing columns (category/date/
#    legislator_name/legislator_id/source_file/index) used on the other
#    sheets (存款/股票/...),
#  - the free-text policy-number note column is dropped in favour of the
#    standard "insurance" property_category value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")
$ws.Activate()

# --- Extend formatting onto the new columns first ----------------------
# Row 1: reuse the existing bold/bordered header style (copied from E1).
$ws.Range("E1").Copy()
$ws.Range("F1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 2 & 3: reuse the existing plain data-row style (copied from E2/E3).
$ws.Range("E2").Copy()
$ws.Range("F2:K2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E3").Copy()
$ws.Range("F3:K3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "date" column holds a plain "2012-05-01" string elsewhere in the
# workbook; force text formatting before assigning so Excel doesn't
# reinterpret it as a date serial number.
$ws.Range("G2:G3").NumberFormat = "@"

# --- Header row (row 1) -------------------------------------------------
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Data row 2 (index 132) ---------------------------------------------
$ws.Range("B2").Value = "中華郵政"
$ws.Range("C2").Value = "6年期步步高升30萬"
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("G2").Value = "2012-05-01"
$ws.Range("H2").Value = "羅淑蕾"
$ws.Range("I2").Value = 1638
$ws.Range("J2").Value = "tmpe6421"
$ws.Range("K2").Value = 132

# --- Data row 3 (index 133) ---------------------------------------------
$ws.Range("B3").Value = "中華郵政"
$ws.Range("E3").Value = "insurance"
$ws.Range("F3").Value = "normal"
$ws.Range("G3").Value = "2012-05-01"
$ws.Range("H3").Value = "羅淑蕾"
$ws.Range("I3").Value = 1638
$ws.Range("J3").Value = "tmpe6421"
$ws.Range("K3").Value = 133
